# Auto-generated edit script: updates market-price derived columns (H-N)
# across multiple crafting-leve profit sheets, per the source diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 148.6
$ws.Range("I8").Value = 76.22221999999999
$ws.Range("K8").Value = 228.66666
$ws.Range("M8").Value = -89.66665999999998
$ws.Range("H18").Value = 1113.8334
$ws.Range("I18").Value = 1170.2
$ws.Range("J18").Value = 832
$ws.Range("K18").Value = 1170.2
$ws.Range("L18").Value = 832
$ws.Range("M18").Value = -886.2
$ws.Range("N18").Value = -1400
$ws.Range("H32").Value = 7844
$ws.Range("J32").Value = 7858.8
$ws.Range("L32").Value = 7858.8
$ws.Range("N32").Value = -8510.799999999999
$ws.Range("H58").Value = 158.5
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").Value = ""
$ws.Range("H111").Value = 1344
$ws.Range("I111").Value = 1342
$ws.Range("K111").Value = 4026
$ws.Range("M111").Value = -959
$ws.Range("H132").Value = 2381.2
$ws.Range("I132").Value = 2632.75
$ws.Range("K132").Value = 7898.25
$ws.Range("M132").Value = -5368.25
$ws.Range("H137").Value = 4857.2856
$ws.Range("I137").Value = 4704.2666
$ws.Range("K137").Value = 14112.7998
$ws.Range("M137").Value = -11562.7998
$ws.Range("H138").Value = 7500.1895
$ws.Range("I138").Value = 4110
$ws.Range("J138").Value = 9151.82
$ws.Range("K138").Value = 12330
$ws.Range("L138").Value = 27455.46
$ws.Range("M138").Value = -7190
$ws.Range("N138").Value = -37735.46

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 670761.4399999999
$ws.Range("I2").Value = 867052.5
$ws.Range("J2").Value = 3371.8
$ws.Range("K2").Value = 867052.5
$ws.Range("L2").Value = 3371.8
$ws.Range("M2").Value = -866939.5
$ws.Range("N2").Value = -3597.8
$ws.Range("H45").Value = 1699.8182
$ws.Range("I45").Value = 1349.8334
$ws.Range("K45").Value = 1349.8334
$ws.Range("M45").Value = -972.8334
$ws.Range("H61").Value = 200001500
$ws.Range("I61").Value = 200001500
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 200001500
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = ""
$ws.Range("N61").Value = -200001288
$ws.Range("H116").Value = 670761.4399999999
$ws.Range("I116").Value = 867052.5
$ws.Range("J116").Value = 3371.8
$ws.Range("K116").Value = 867052.5
$ws.Range("L116").Value = 3371.8
$ws.Range("M116").Value = -864758.5
$ws.Range("N116").Value = -7959.8
$ws.Range("H132").Value = 4225478.5
$ws.Range("I132").Value = 2384212.5
$ws.Range("K132").Value = 7152637.5
$ws.Range("M132").Value = -7150107.5
$ws.Range("H136").Value = 200001500
$ws.Range("I136").Value = 200001500
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 600004500
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = ""
$ws.Range("N136").Value = -600001950

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 670761.4399999999
$ws.Range("I3").Value = 867052.5
$ws.Range("J3").Value = 3371.8
$ws.Range("K3").Value = 867052.5
$ws.Range("L3").Value = 3371.8
$ws.Range("M3").Value = -866938.5
$ws.Range("N3").Value = -3599.8
$ws.Range("H86").Value = 2533.0952
$ws.Range("I86").Value = 2615.3845
$ws.Range("J86").Value = 2399.375
$ws.Range("K86").Value = 2615.3845
$ws.Range("L86").Value = 2399.375
$ws.Range("M86").Value = -1492.3845
$ws.Range("N86").Value = -4645.375
$ws.Range("H89").Value = 2533.0952
$ws.Range("I89").Value = 2615.3845
$ws.Range("J89").Value = 2399.375
$ws.Range("K89").Value = 13076.9225
$ws.Range("L89").Value = 11996.875
$ws.Range("M89").Value = -7460.922500000001
$ws.Range("N89").Value = -23228.875
$ws.Range("H105").Value = 3224.7778
$ws.Range("I105").Value = 2662.3635
$ws.Range("K105").Value = 2662.3635
$ws.Range("M105").Value = -915.3634999999999
$ws.Range("H134").Value = 17595406
$ws.Range("I134").Value = 18223746
$ws.Range("J134").Value = 1899
$ws.Range("K134").Value = 54671238
$ws.Range("L134").Value = 5697
$ws.Range("M134").Value = -54668703
$ws.Range("N134").Value = -10767

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 42333
$ws.Range("J28").Value = 42333
$ws.Range("L28").Value = 42333
$ws.Range("N28").Value = -42823
$ws.Range("H31").Value = 3733.9365
$ws.Range("I31").Value = 2350.7222
$ws.Range("J31").Value = 12033.223
$ws.Range("K31").Value = 2350.7222
$ws.Range("L31").Value = 12033.223
$ws.Range("M31").Value = -2055.7222
$ws.Range("N31").Value = -12623.223
$ws.Range("H34").Value = 3733.9365
$ws.Range("I34").Value = 2350.7222
$ws.Range("J34").Value = 12033.223
$ws.Range("K34").Value = 2350.7222
$ws.Range("L34").Value = 12033.223
$ws.Range("M34").Value = -2148.7222
$ws.Range("N34").Value = -12437.223
$ws.Range("H122").Value = 3128.3333
$ws.Range("I122").Value = 2988.7368
$ws.Range("K122").Value = 8966.2104
$ws.Range("M122").Value = -6516.2104

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 500
$ws.Range("I58").Value = 500
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1500
$ws.Range("L58").Value = ""
$ws.Range("M58").Value = -1372
$ws.Range("N58").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3259.5454
$ws.Range("J22").Value = 2916.2
$ws.Range("L22").Value = 2916.2
$ws.Range("N22").Value = -3506.2
$ws.Range("H27").Value = 3259.5454
$ws.Range("J27").Value = 2916.2
$ws.Range("L27").Value = 2916.2
$ws.Range("N27").Value = -3130.2
$ws.Range("H46").Value = 1665.1666
$ws.Range("I46").Value = 1674.25
$ws.Range("K46").Value = 1674.25
$ws.Range("M46").Value = -1486.25
$ws.Range("H55").Value = 539.1
$ws.Range("J55").Value = 921
$ws.Range("L55").Value = 921
$ws.Range("N55").Value = -1267
$ws.Range("H100").Value = 22178078
$ws.Range("I100").Value = 28514180
$ws.Range("K100").Value = 28514180
$ws.Range("M100").Value = -28513639
$ws.Range("H136").Value = 2236.0625
$ws.Range("I136").Value = 2190.05
$ws.Range("J136").Value = 2312.75
$ws.Range("K136").Value = 6570.150000000001
$ws.Range("L136").Value = 6938.25
$ws.Range("M136").Value = -4020.150000000001
$ws.Range("N136").Value = -12038.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1696.625
$ws.Range("J107").Value = 3476
$ws.Range("L107").Value = 10428
$ws.Range("N107").Value = -14268
$ws.Range("H113").Value = 625.5454999999999
$ws.Range("I113").Value = 633.9048
$ws.Range("J113").Value = 450
$ws.Range("K113").Value = 1901.7144
$ws.Range("L113").Value = 1350
$ws.Range("M113").Value = 268.2855999999999
$ws.Range("N113").Value = -5690

Write-Host "All cell updates applied."